$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 5468.5
$ws.Range("I76").Value = 5466.4
$ws.Range("J76").Value = 5500
$ws.Range("K76").Value = 5466.4
$ws.Range("L76").Value = 5500
$ws.Range("M76").Value = -5151.4
$ws.Range("N76").Value = -6130
# Row 79
$ws.Range("H79").Value = 5468.5
$ws.Range("I79").Value = 5466.4
$ws.Range("J79").Value = 5500
$ws.Range("K79").Value = 5466.4
$ws.Range("L79").Value = 5500
$ws.Range("M79").Value = -4374.4
$ws.Range("N79").Value = -7684
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 106
$ws.Range("H106").Value = 2771
$ws.Range("I106").Value = 2771
$ws.Range("K106").Value = 2771
$ws.Range("M106").Value = -2140
# Row 112
$ws.Range("H112").Value = 1148.8334
$ws.Range("J112").Value = 1133.5652
$ws.Range("L112").Value = 3400.6956
$ws.Range("N112").Value = -5616.6956
# Row 113
$ws.Range("H113").Value = 3506.2856
$ws.Range("I113").Value = 2899.2222
$ws.Range("J113").Value = 4599
$ws.Range("K113").Value = 2899.2222
$ws.Range("L113").Value = 4599
$ws.Range("M113").Value = 354.7777999999998
$ws.Range("N113").Value = -11107
# Row 138
$ws.Range("H138").Value = 3662.94
$ws.Range("I138").Value = 2963.1428
$ws.Range("J138").Value = 3935.0833
$ws.Range("K138").Value = 8889.428400000001
$ws.Range("L138").Value = 11805.2499
$ws.Range("M138").Value = -3749.428400000001
$ws.Range("N138").Value = -22085.2499

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3982.195
$ws.Range("I32").Value = 3270.7104
$ws.Range("K32").Value = 3270.7104
$ws.Range("M32").Value = -2983.7104
# Row 61
$ws.Range("H61").Value = 5992.143
$ws.Range("I61").Value = 2579.8
$ws.Range("K61").Value = 2579.8
$ws.Range("M61").Value = -2367.8
# Row 74
$ws.Range("H74").Value = 2678.1738
$ws.Range("I74").Value = 2454.7778
$ws.Range("K74").Value = 2454.7778
$ws.Range("M74").Value = -1580.7778
# Row 77
$ws.Range("H77").Value = 2678.1738
$ws.Range("I77").Value = 2454.7778
$ws.Range("K77").Value = 12273.889
$ws.Range("M77").Value = -7905.888999999999
# Row 135
$ws.Range("H135").Value = 307499.5
$ws.Range("J135").Value = 307499.5
$ws.Range("L135").Value = 307499.5
$ws.Range("N135").Value = -317639.5
# Row 136
$ws.Range("H136").Value = 5992.143
$ws.Range("I136").Value = 2579.8
$ws.Range("K136").Value = 7739.400000000001
$ws.Range("M136").Value = -5189.400000000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 427.54544
$ws.Range("I94").Value = 430.3
$ws.Range("K94").Value = 430.3
$ws.Range("M94").Value = 20.69999999999999
# Row 105
$ws.Range("H105").Value = 2705.7827
$ws.Range("I105").Value = 2626.5557
$ws.Range("K105").Value = 2626.5557
$ws.Range("M105").Value = -879.5556999999999
# Row 107
$ws.Range("H107").Value = 3511.5557
$ws.Range("I107").Value = 2913
$ws.Range("J107").Value = 5606.5
$ws.Range("K107").Value = 2913
$ws.Range("L107").Value = 5606.5
$ws.Range("M107").Value = -993
$ws.Range("N107").Value = -9446.5
# Row 132
$ws.Range("H132").Value = 84499.5
$ws.Range("J132").Value = 84499.5
$ws.Range("L132").Value = 84499.5
$ws.Range("N132").Value = -94619.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 3999.5
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
# Row 19
$ws.Range("H19").Value = 158.75
$ws.Range("I19").Value = 158.75
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 158.75
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 11.25
$ws.Range("N19").ClearContents()
# Row 24
$ws.Range("H24").Value = 158.75
$ws.Range("I24").Value = 158.75
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 158.75
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 11.25
$ws.Range("N24").ClearContents()
# Row 31
$ws.Range("H31").Value = 6852.143
$ws.Range("I31").Value = 5378.8
$ws.Range("J31").Value = 7670.6665
$ws.Range("K31").Value = 5378.8
$ws.Range("L31").Value = 7670.6665
$ws.Range("M31").Value = -5083.8
$ws.Range("N31").Value = -8260.666499999999
# Row 34
$ws.Range("H34").Value = 6852.143
$ws.Range("I34").Value = 5378.8
$ws.Range("J34").Value = 7670.6665
$ws.Range("K34").Value = 5378.8
$ws.Range("L34").Value = 7670.6665
$ws.Range("M34").Value = -5176.8
$ws.Range("N34").Value = -8074.6665
# Row 122
$ws.Range("H122").Value = 1143.6666
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 47
$ws.Range("H47").Value = 1555.5555
$ws.Range("I47").Value = 1500
$ws.Range("J47").Value = 2000
$ws.Range("K47").Value = 4500
$ws.Range("L47").Value = 6000
$ws.Range("M47").Value = -4069
$ws.Range("N47").Value = -6862
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 112
$ws.Range("H112").Value = 10131.571
$ws.Range("J112").Value = 19998.334
$ws.Range("L112").Value = 59995.00199999999
$ws.Range("N112").Value = -62211.00199999999
# Row 113
$ws.Range("H113").Value = 1521.3889
$ws.Range("J113").Value = 1364.0714
$ws.Range("L113").Value = 4092.2142
$ws.Range("N113").Value = -8432.2142

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 103957.14
$ws.Range("I5").Value = 103957.14
$ws.Range("K5").Value = 103957.14
$ws.Range("M5").Value = -103845.14
# Row 80
$ws.Range("H80").Value = 13295.2
$ws.Range("I80").Value = 2972
$ws.Range("J80").Value = 17719.428
$ws.Range("K80").Value = 2972
$ws.Range("L80").Value = 17719.428
$ws.Range("M80").Value = -1974
$ws.Range("N80").Value = -19715.428
# Row 83
$ws.Range("H83").Value = 13295.2
$ws.Range("I83").Value = 2972
$ws.Range("J83").Value = 17719.428
$ws.Range("K83").Value = 14860
$ws.Range("L83").Value = 88597.14
$ws.Range("M83").Value = -9868
$ws.Range("N83").Value = -98581.14
# Row 97
$ws.Range("H97").Value = 726.8889
$ws.Range("J97").Value = 749.6667
$ws.Range("L97").Value = 749.6667
$ws.Range("N97").Value = -1741.6667
# Row 107
$ws.Range("H107").Value = 329.125
$ws.Range("I107").Value = 233.28572
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 233.28572
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1686.71428
$ws.Range("N107").Value = -4840

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 1196973.6
$ws.Range("I2").Value = 143971
$ws.Range("J2").Value = 2074475.9
$ws.Range("K2").Value = 143971
$ws.Range("L2").Value = 2074475.9
$ws.Range("M2").Value = -143859
$ws.Range("N2").Value = -2074699.9
# Row 40
$ws.Range("H40").Value = 5000.875
$ws.Range("I40").Value = 4502.3335
$ws.Range("K40").Value = 4502.3335
$ws.Range("M40").Value = -4366.3335
# Row 61
$ws.Range("H61").Value = 6955.231
$ws.Range("I61").Value = 5754.4287
$ws.Range("J61").Value = 11998.6
$ws.Range("K61").Value = 5754.4287
$ws.Range("L61").Value = 11998.6
$ws.Range("M61").Value = -5552.4287
$ws.Range("N61").Value = -12402.6
# Row 113
$ws.Range("H113").Value = 6955.231
$ws.Range("I113").Value = 5754.4287
$ws.Range("J113").Value = 11998.6
$ws.Range("K113").Value = 5754.4287
$ws.Range("L113").Value = 11998.6
$ws.Range("M113").Value = -3584.4287
$ws.Range("N113").Value = -16338.6
# Row 132
$ws.Range("H132").Value = 4258.25
$ws.Range("I132").Value = 4011.3333
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 12033.9999
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -9503.999899999999
$ws.Range("N132").Value = -20057

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 11
$ws.Range("H11").Value = 255000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
# Row 113
$ws.Range("H113").Value = 294.8
$ws.Range("I113").Value = 118.5
$ws.Range("K113").Value = 355.5
$ws.Range("M113").Value = 1814.5
# Row 122
$ws.Range("H122").Value = 2037.8334
$ws.Range("I122").Value = 2146.4
$ws.Range("K122").Value = 6439.200000000001
$ws.Range("M122").Value = -3989.200000000001
# Row 132
$ws.Range("H132").Value = 9096.5
$ws.Range("I132").Value = 7128.6665
$ws.Range("K132").Value = 21385.9995
$ws.Range("M132").Value = -18855.9995
# Row 136
$ws.Range("H136").Value = 2555.4546
$ws.Range("I136").Value = 2555.4546
$ws.Range("K136").Value = 7666.3638
$ws.Range("M136").Value = -5116.3638

